# Nalco aluminium ingot price sheet: a new day's price row is published
# (03-12-2025), so every existing row shifts down by one and a fresh
# row 2 is written with that day's figures (same price/circular info as
# the previous top row, only the quoted "Date" column moves on).
# The row that falls off the bottom of the original range (old row 119)
# is preserved as the new row 120, complete with its own hyperlink.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the values currently sitting in row 2 - these become the
# "Circular Date"/"Circular Link"/etc. for the freshly inserted row.
$oldB2 = $ws.Range("B2").Value()
$oldC2 = $ws.Range("C2").Value()
$oldD2 = $ws.Range("D2").Value()
$oldE2 = $ws.Range("E2").Value()
$oldF2 = $ws.Range("F2").Value()

# Push everything down one row. This alone grows the sheet from
# A1:F119 to A1:F120 and relocates all existing hyperlinks correctly.
$ws.Rows.Item(2).Insert()

# Force the date-like text in column A/E to stay plain text instead of
# Excel auto-parsing it into a date serial (matches the rest of the
# column, which stores dates as literal strings).
$ws.Range("A2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"

$ws.Range("A2").Value = "03-12-2025"
$ws.Range("B2").Value = $oldB2
$ws.Range("C2").Value = $oldC2
$ws.Range("D2").Value = $oldD2
$ws.Range("E2").Value = $oldE2
$ws.Range("F2").Value = $oldF2

# Re-apply the formatting used by the rest of the data rows (the insert
# plus the NumberFormat tweak above left row 2 with ad-hoc styling).
$ws.Range("A3:F3").Copy()
$ws.Range("A2:F2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# The row that used to be last (07-08-2025) is now row 120; it needs
# its own hyperlink on column F, same as every other data row.
$lastRow = 120
$lastUrl = $ws.Range("F" + $lastRow).Value()
$ws.Hyperlinks.Add($ws.Range("F" + $lastRow), $lastUrl)

# Adding the hyperlink re-styles the cell with the built-in Hyperlink
# look; bring it back in line with the rest of the column.
$ws.Range("F119").Copy()
$ws.Range("F" + $lastRow).PasteSpecial(-4122)
$excel.CutCopyMode = $false
